$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns: email, company, comment
$ws.Range("C1").Value = "email"
$ws.Range("D1").Value = "company"
$ws.Range("E1").Value = "comment"
# Match the bold header formatting already used by A1/B1
$ws.Range("C1:E1").Font.Bold = $true

# Row 2 data for the new columns
$ws.Range("A2").Style = "Normal"
$ws.Range("C2").Value = "hello@gmail.com"
$ws.Range("D2").Value = "Network"
$ws.Range("E2").Value = "Hello I am Hello User"

# Turn the email into a mailto hyperlink, then restore the Hyperlink style
# on the cell (Hyperlinks.Add re-applies its own formatting, so we set the
# named style afterwards to land on the workbook's existing Hyperlink xf).
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:hello@gmail.com")
$ws.Range("C2").Style = "Hyperlink"

# Move the active selection like the author's last action in the sheet
$ws.Range("C10").Select() | Out-Null
